# Insert a new data row at row 34 (this shifts the existing rows 34..123 down to 35..124,
# which matches the rest of the diff where every row's values equal the prior row's values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("34:34").Insert()

# Fill in the new row 34 with the newly-added record.
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 44804
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100106
$ws.Range("H34").Value = "Oleaginosos"
$ws.Range("I34").Value = 100106002
$ws.Range("J34").Value = "Palta"
$ws.Range("K34").Value = "Fuerte"
$ws.Range("L34").Value = "Tercera"
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 50000
$ws.Range("O34").Value = 52000
$ws.Range("P34").Value = 51000
$ws.Range("Q34").Value = "`$/caja 25 kilos"
$ws.Range("R34").Value = "Región de Coquimbo"
$ws.Range("S34").Value = 2040
$ws.Range("T34").Value = 25
